# Weekly update: insert a new data row at row 93 (most recent week's price
# record for "Zapallo italiano" at Terminal La Palmera de La Serena),
# pushing the previously existing rows 93:195 down to 94:196.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 93, shifting rows 93:195 -> 94:196
$ws.Rows("93:93").Insert()

# Populate the newly inserted row 93 with the new weekly record
$ws.Range("A93").Value = 8
$ws.Range("B93").Value = "Terminal La Palmera de La Serena"
$ws.Range("C93").Value = "Coquimbo"
$ws.Range("D93").Value = 44494
$ws.Range("E93").Value = 4
$ws.Range("F93").Value = 100112032
$ws.Range("G93").Value = "Zapallo italiano"
$ws.Range("H93").Value = "Sin especificar"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 500
$ws.Range("K93").Value = 12000
$ws.Range("L93").Value = 12500
$ws.Range("M93").Value = 12250
$ws.Range("N93").Value = "`$/caja 70 unidades"
$ws.Range("O93").Value = "Provincia de Limarí"
$ws.Range("P93").Value = 175
$ws.Range("Q93").Value = 70
$ws.Range("R93").Value = "Hortaliza"

# Apply the same date style (yyyy-mm-dd hh:mm:ss, style used by column D)
# as the neighboring date cells to the newly inserted D93
$ws.Range("D93").NumberFormat = $ws.Range("D94").NumberFormat
